$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition list)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 4027
$ws1.Range("F4").Value = 2357
$ws1.Range("F5").Value = 469
$ws1.Range("F9").Value = 197
$ws1.Range("F11").Value = 69
$ws1.Range("F13").Value = 1499
$ws1.Range("F15").Value = 2832

# Sheet "全部类型" (all types list) - same events, different row offsets
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 4027
$ws4.Range("F4").Value = 2357
$ws4.Range("F5").Value = 469
$ws4.Range("F10").Value = 197
$ws4.Range("F12").Value = 69
$ws4.Range("F16").Value = 1499
$ws4.Range("F18").Value = 2832
